$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hybris")
$ws1.Columns.Item(2).ColumnWidth = 19.5703125
$ws1.Activate()
$ws1.Range("C17").Select()
